$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Logs")

# Append the new log row (row 22) mirroring the previous "Kun jij dit even
# regelen?" entry, but with its own reply timestamp.
$ws.Range("A22").Value = "Kun jij dit even regelen?"
$ws.Range("B22").Value = "mailmind.test@zohomail.eu"
$ws.Range("C22").Value = "Testmail #1: Kun jij dit even regelen?"
$ws.Range("D22").Value = "Planning / Afspraak"
$ws.Range("E22").Value = "Bedankt, we hebben dit doorgestuurd naar planning@bedrijf.nl."
$ws.Range("F22").Value = "2025-08-05 18:06:39"
$ws.Range("G22").Value = "Ja"
$ws.Range("H22").Value = "Ja"
$ws.Range("I22").Value = "Nee"
$ws.Range("J22").Value = "Nee"

# Extend the conditional formatting ranges so the new row is covered too.
$cols = @("D", "G", "H", "I", "J")
foreach ($col in $cols) {
    $oldRange = $col + "2:" + $col + "21"
    $newRange = $col + "2:" + $col + "22"
    $fc = $ws.Range($oldRange).FormatConditions.Item(1)
    $fc.ModifyAppliesToRange($ws.Range($newRange))
}

# Update the Dashboard summary count for "Planning / Afspraak".
$dash = $wb.Worksheets.Item("Dashboard")
$dash.Range("B2").Value = 16
